# Update the "fix cost" results row (row 2) on each yearly sheet with the
# latest values received from the server.

$wb = $excel.ActiveWorkbook

$sheetNames = @("2025", "2030", "2035", "2040", "2045", "2050")

$sheetUpdates = @(
    @{
        "A2" = 0
        "B2" = 135.9602344995041
        "E2" = 346.7127090767025
        "G2" = 182.1583285348912
        "I2" = 777.8494427183999
        "L2" = 0
        "M2" = 0
        "N2" = 125.3734990721661
        "O2" = 183.6538381306478
    },
    @{
        "A2" = 101.9170951152791
        "B2" = 362.9066076366151
        "E2" = 1313.42977404479
        "G2" = 182.1583285348912
        "I2" = 2227.614065720369
        "L2" = 0
        "M2" = 0
        "N2" = 337.5357088227005
        "O2" = 363.7792159379853
    },
    @{
        "A2" = 234.6760704606383
        "B2" = 714.3598950717846
        "E2" = 2767.358511932531
        "G2" = 182.1583285348912
        "I2" = 4377.552022949434
        "L2" = 0
        "M2" = 0
        "N2" = 735.9252380194357
        "O2" = 664.3811208296214
    },
    @{
        "A2" = 234.6760704606383
        "B2" = 714.3598950717846
        "E2" = 2767.358511932531
        "G2" = 182.1583285348912
        "I2" = 4377.552022949434
        "L2" = 0
        "M2" = 0
        "N2" = 735.9252380194357
        "O2" = 684.6978974915725
    },
    @{
        "A2" = 292.8744314431498
        "B2" = 714.3598950717846
        "E2" = 2767.358511932531
        "G2" = 182.1583285348912
        "I2" = 4377.552022949434
        "L2" = 0
        "M2" = 0
        "N2" = 735.9252380194357
        "O2" = 701.6984641353358
    },
    @{
        "A2" = 292.8744314431498
        "B2" = 714.3598950717846
        "E2" = 2767.358511932531
        "G2" = 182.1583285348912
        "I2" = 4377.552022949434
        "L2" = 0
        "M2" = 0
        "N2" = 735.9252380194357
        "O2" = 701.6984641353358
    }
)

for ($i = 0; $i -lt $sheetNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($sheetNames[$i])
    $cellValues = $sheetUpdates[$i]
    foreach ($entry in $cellValues.GetEnumerator()) {
        $ws.Range($entry.Key).Value = $entry.Value
    }
}
